$d = $word.ActiveDocument

# 1) Fix typo: "2.5 c.c. d'huile" -> "2.5 c.s. d'huile" (it's a tablespoon
#    of oil, like the other quantities in that row, not a teaspoon).
#    Locate the specific table cell by its content so the other "c.c."
#    occurrences elsewhere in the table (salt, yeast) are left untouched.
foreach ($tbl in $d.Tables) {
    foreach ($cell in $tbl.Range.Cells) {
        if ($cell.Range.Text -like "*c.c. d*huile*") {
            $cell.Range.Find.Execute(" c.c. d", $true, $false, $false, $false,
                                      $false, $true, 1, $false, " c.s. d", 1)
        }
    }
}

# 2) Remove the extra "Attention: On ne peut PAS remplacer l'eau et le lait
#    en poudre par 300ml de lait normal." paragraph entirely (the whole
#    paragraph, including its trailing paragraph mark).
foreach ($p in $d.Content.Paragraphs) {
    if ($p.Range.Text -like "*300ml de lait normal*") {
        $p.Range.Delete()
        break
    }
}
